$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 10; this shifts the existing rows 10-13 down to 11-14,
# keeping all of their data intact.
$ws.Rows.Item(10).Insert()

# Fill the new row 10 with the new reading (same template as the old row 10, new date/volume).
$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(10, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(10, 4).Value = 44460
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 5).Value = 15
$ws.Cells.Item(10, 6).Value = 100112003
$ws.Cells.Item(10, 7).Value = "Ajo"
$ws.Cells.Item(10, 8).Value = "Chino"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 300
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 16000
$ws.Cells.Item(10, 13).Value = 15500
$ws.Cells.Item(10, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(10, 15).Value = "China"
$ws.Cells.Item(10, 16).Value = 1550
$ws.Cells.Item(10, 17).Value = 10
$ws.Cells.Item(10, 18).Value = "Hortaliza"
